$d = $word.ActiveDocument

$d.Content.Find.Execute("855÷2=427, 1", $true, $false, $false, $false, $false, $true, 1, $false, "778÷9=86, 4", 2) | Out-Null
$d.Content.Find.Execute("739÷9=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "413÷2=206, 1", 2) | Out-Null
$d.Content.Find.Execute("140÷5=28, 0", $true, $false, $false, $false, $false, $true, 1, $false, "754÷3=251, 1", 2) | Out-Null
$d.Content.Find.Execute("681÷9=75, 6", $true, $false, $false, $false, $false, $true, 1, $false, "843÷8=105, 3", 2) | Out-Null
$d.Content.Find.Execute("896÷6=149, 2", $true, $false, $false, $false, $false, $true, 1, $false, "389÷3=129, 2", 2) | Out-Null
$d.Content.Find.Execute("101÷6=16, 5", $true, $false, $false, $false, $false, $true, 1, $false, "684÷7=97, 5", 2) | Out-Null
$d.Content.Find.Execute("213÷9=23, 6", $true, $false, $false, $false, $false, $true, 1, $false, "259÷9=28, 7", 2) | Out-Null
$d.Content.Find.Execute("107÷7=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "825÷5=165, 0", 2) | Out-Null
$d.Content.Find.Execute("176÷6=29, 2", $true, $false, $false, $false, $false, $true, 1, $false, "905÷8=113, 1", 2) | Out-Null
$d.Content.Find.Execute("651÷9=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "841÷8=105, 1", 2) | Out-Null
$d.Content.Find.Execute("660÷3=220, 0", $true, $false, $false, $false, $false, $true, 1, $false, "199÷8=24, 7", 2) | Out-Null
$d.Content.Find.Execute("349÷2=174, 1", $true, $false, $false, $false, $false, $true, 1, $false, "319÷4=79, 3", 2) | Out-Null
$d.Content.Find.Execute("332÷5=66, 2", $true, $false, $false, $false, $false, $true, 1, $false, "712÷7=101, 5", 2) | Out-Null
$d.Content.Find.Execute("618÷8=77, 2", $true, $false, $false, $false, $false, $true, 1, $false, "991÷6=165, 1", 2) | Out-Null
$d.Content.Find.Execute("104÷2=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "372÷3=124, 0", 2) | Out-Null
$d.Content.Find.Execute("450÷2=225, 0", $true, $false, $false, $false, $false, $true, 1, $false, "188÷8=23, 4", 2) | Out-Null
$d.Content.Find.Execute("872÷4=218, 0", $true, $false, $false, $false, $false, $true, 1, $false, "127÷2=63, 1", 2) | Out-Null
$d.Content.Find.Execute("189÷6=31, 3", $true, $false, $false, $false, $false, $true, 1, $false, "543÷9=60, 3", 2) | Out-Null
$d.Content.Find.Execute("104÷9=11, 5", $true, $false, $false, $false, $false, $true, 1, $false, "214÷4=53, 2", 2) | Out-Null
$d.Content.Find.Execute("713÷3=237, 2", $true, $false, $false, $false, $false, $true, 1, $false, "421÷6=70, 1", 2) | Out-Null
$d.Content.Find.Execute("679÷8=84, 7", $true, $false, $false, $false, $false, $true, 1, $false, "395÷8=49, 3", 2) | Out-Null
$d.Content.Find.Execute("277÷4=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "233÷3=77, 2", 2) | Out-Null
$d.Content.Find.Execute("139÷3=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "743÷3=247, 2", 2) | Out-Null
$d.Content.Find.Execute("493÷4=123, 1", $true, $false, $false, $false, $false, $true, 1, $false, "138÷2=69, 0", 2) | Out-Null
$d.Content.Find.Execute("555÷4=138, 3", $true, $false, $false, $false, $false, $true, 1, $false, "145÷6=24, 1", 2) | Out-Null
